$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40: new work log entry (copy date format from row 39, then set value)
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = 41702
$ws.Range("B40").Value = "Team meeting. Recorded meeting minutes. Created formal intermedial designs and team organization diagram from meeting notes."
$ws.Range("I40").Value = 2.5

# Row 41: new work log entry
$ws.Range("A39").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A41").Value = 41704
$ws.Range("B41").Value = "Team meeting and client meeting. Recorded meeting minutes. Drafted architectural design document. Participated in initial discussion of detailed design challenges."
$ws.Range("I41").Value = 3.25

# Update the view state to match scrolled position / selection
$ws.Range("B43:H43").Select()
